# This script applies the latest cryptocurrency market-data snapshot to the
# "cryptos" worksheet (price + 1h volume change per coin), including the two
# row swaps (Fetch.AI/Dai at rows 27-28, and Bittensor/USDe at rows 48-49)
# that resulted from the underlying ranking changing order.
#
# Price/percentage values are written as literal text (matching the source
# data, which stores them as plain strings such as "1.00" or "0.0000105")
# rather than as numbers, so formatting/precision is preserved exactly.
# NumberFormat is temporarily switched to text ("@") while assigning the
# value so Excel does not auto-convert numeric-looking text into a float,
# then the cell's original Style is restored so no visible formatting
# change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '61.740.24'
Set-TextValue $ws 'E2' '  -1.30%  '

Set-TextValue $ws 'D3' '2.906.70'
Set-TextValue $ws 'E3' '  -1.73%  '

Set-TextValue $ws 'D4' '1.00'
Set-TextValue $ws 'E4' '  -0.08%  '

Set-TextValue $ws 'D5' '586.44'
Set-TextValue $ws 'E5' '  -1.39%  '

Set-TextValue $ws 'D6' '146.19'
Set-TextValue $ws 'E6' '  +1.02%  '

Set-TextValue $ws 'E7' '  -0.07%  '

Set-TextValue $ws 'D8' '0.507'
Set-TextValue $ws 'E8' '  +0.15%  '

Set-TextValue $ws 'D9' '2.906.62'
Set-TextValue $ws 'E9' '  -1.67%  '

Set-TextValue $ws 'D10' '6.92'
Set-TextValue $ws 'E10' '  -5.61%  '

Set-TextValue $ws 'D11' '0.149'
Set-TextValue $ws 'E11' '  +4.43%  '

Set-TextValue $ws 'D12' '0.434'
Set-TextValue $ws 'E12' '  -2.71%  '

Set-TextValue $ws 'E13' '  +1.62%  '

Set-TextValue $ws 'D14' '32.81'
Set-TextValue $ws 'E14' '  -1.45%  '

Set-TextValue $ws 'E15' '  -0.90%  '

Set-TextValue $ws 'D16' '3.388.36'
Set-TextValue $ws 'E16' '  -1.86%  '

Set-TextValue $ws 'D17' '61.797.67'
Set-TextValue $ws 'E17' '  -1.16%  '

Set-TextValue $ws 'D18' '6.62'
Set-TextValue $ws 'E18' '  -1.35%  '

Set-TextValue $ws 'D19' '2.906.27'
Set-TextValue $ws 'E19' '  -1.88%  '

Set-TextValue $ws 'D20' '435.54'
Set-TextValue $ws 'E20' '  -1.32%  '

Set-TextValue $ws 'D21' '13.37'
Set-TextValue $ws 'E21' '  -0.31%  '

Set-TextValue $ws 'E22' '  -1.59%  '

Set-TextValue $ws 'D23' '6.94'
Set-TextValue $ws 'E23' '  -2.17%  '

Set-TextValue $ws 'D24' '81.09'
Set-TextValue $ws 'E24' '  -0.87%  '

Set-TextValue $ws 'D25' '11.95'
Set-TextValue $ws 'E25' '  -0.20%  '

Set-TextValue $ws 'D26' '10.31'
Set-TextValue $ws 'E26' '  -7.20%  '

Set-TextValue $ws 'B27' 'Fetch.AI'
Set-TextValue $ws 'C27' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws 'D27' '2.08'
Set-TextValue $ws 'E27' '  -2.28%  '

Set-TextValue $ws 'B28' 'Dai'
Set-TextValue $ws 'C28' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D28' '0.999'
Set-TextValue $ws 'E28' '  -0.13%  '

Set-TextValue $ws 'D29' '0.0000105'
Set-TextValue $ws 'E29' '  +20.41%  '

Set-TextValue $ws 'D30' '7.15'
Set-TextValue $ws 'E30' '  +1.56%  '

Set-TextValue $ws 'D31' '2.56'
Set-TextValue $ws 'E31' '  -1.93%  '

Set-TextValue $ws 'E32' '  -0.68%  '

Set-TextValue $ws 'E33' '  +1.16%  '

Set-TextValue $ws 'D34' '0.999'
Set-TextValue $ws 'E34' '  -0.25%  '

Set-TextValue $ws 'D35' '25.90'
Set-TextValue $ws 'E35' '  -2.87%  '

Set-TextValue $ws 'E36' '  -1.98%  '

Set-TextValue $ws 'E37' '  -2.10%  '

Set-TextValue $ws 'E38' '  +3.62%  '

Set-TextValue $ws 'E39' '  -1.00%  '

Set-TextValue $ws 'E40' '  -2.65%  '

Set-TextValue $ws 'D41' '8.37'
Set-TextValue $ws 'E41' '  -2.09%  '

Set-TextValue $ws 'E42' '  -0.92%  '

Set-TextValue $ws 'E43' '  -2.88%  '

Set-TextValue $ws 'D44' '39.02'
Set-TextValue $ws 'E44' '  +0.18%  '

Set-TextValue $ws 'D45' '2.698.98'
Set-TextValue $ws 'E45' '  -0.72%  '

Set-TextValue $ws 'D46' '133.76'
Set-TextValue $ws 'E46' '  -1.06%  '

Set-TextValue $ws 'D47' '0.0336'
Set-TextValue $ws 'E47' '  -1.63%  '

Set-TextValue $ws 'B48' 'USDe'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws 'D48' '1.00'
Set-TextValue $ws 'E48' '  +0.01%  '

Set-TextValue $ws 'B49' 'Bittensor'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws 'D49' '341.43'
Set-TextValue $ws 'E49' '  -6.25%  '

Set-TextValue $ws 'E50' '  -1.81%  '

Set-TextValue $ws 'D51' '22.27'
Set-TextValue $ws 'E51' '  -2.69%  '
